# Generate Report for Handoff
#
# Updates the localization-status report:
#  - the previously-failing doc (b868da23-...) is replaced by a newly
#    generated doc (9cc9d2d9-...) that is now "Ready for handoff"
#  - a brand new source file (ffff877fdf4d-...) shows up, also ready for
#    handoff, with freshly generated xlf hand-off artifacts
#  - the ".localization-config" bookkeeping row shifts down to make room
#    for the new row
#
# Applied identically to the Overview sheet (file-name x locale summary)
# and to each per-locale detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$repoBlob = "https://github.com/OpenLocalizationTest/oltest/blob"
$docCommit = "ff9ba2c73e2a4380023404c65e3a8373d92104bb"
$cfgCommit = "f379dd7618b7716233029521916f0754331bcc64"

$newDocMd    = "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md"
$newCfgMd    = "ffff877fdf4d-094b-49f3-861a-ea9305a4a539.md"
$localConfig = ".localization-config"

$statusReady  = "Ready for handoff"
$statusNotLoc = "Not to be localized"
$statusInclude = "Include"
$statusIgnored = "Ignored"

$dtZero = "0001-01-01 00:00:00"

$xlfBase = "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Clear out the old hyperlinks up front so re-adding them below does not
# leave stale/duplicate <hyperlink> entries behind.
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newDocMd
$ws.Range("B2").Value = $statusReady
$ws.Range("C2").Value = $statusReady

$ws.Range("A3").Value = $newCfgMd
$ws.Range("B3").Value = $statusReady
$ws.Range("C3").Value = $statusReady

$ws.Range("A4").Value = $localConfig
$ws.Range("B4").Value = $statusNotLoc
$ws.Range("C4").Value = $statusNotLoc

$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBlob/$docCommit/e2e/$newDocMd", "", "", $newDocMd)
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBlob/$docCommit/e2e/$newCfgMd", "", "", $newCfgMd)
$ws.Hyperlinks.Add($ws.Range("A4"), "$repoBlob/$cfgCommit/$localConfig", "", "", $localConfig)

# ---------------------------------------------------------------------
# Helper sets for the two detail sheets (zh-cn / de-de): only the xlf
# file name and handoff datetime differ between locales.
# ---------------------------------------------------------------------
$locales = @(
    @{ SheetName = "zh-cn"; XlfDisplay = "$xlfBase.zh-cn.xlf"; Datetime = "2016-02-18 04:03:46" },
    @{ SheetName = "de-de"; XlfDisplay = "$xlfBase.de-de.xlf"; Datetime = "2016-02-18 04:03:58" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.SheetName)
    $xlfDisplay = $locale.XlfDisplay
    $dt = $locale.Datetime

    # Clear out the old hyperlinks up front (same reasoning as above).
    $ws.Hyperlinks.Delete()

    # Row 2: the doc that previously failed handoff is now ready.
    $ws.Range("A2").Value = $newDocMd
    $ws.Range("B2").Value = $statusReady
    $ws.Range("C2").Value = $xlfDisplay
    $ws.Range("D2").Value = $dt
    $ws.Range("G2").Value = $dtZero
    $ws.Range("H2").Value = $statusInclude

    # Row 3: newly discovered source file, also ready for handoff.
    $ws.Range("A3").Value = $newCfgMd
    $ws.Range("B3").Value = $statusReady
    $ws.Range("C3").Value = $xlfDisplay
    $ws.Range("D3").Value = $dt
    $ws.Range("G3").Value = $dtZero
    $ws.Range("H3").Value = $statusInclude

    # Row 4: the ".localization-config" bookkeeping row, pushed down.
    $ws.Range("A4").Value = $localConfig
    $ws.Range("B4").Value = $statusNotLoc
    $ws.Range("D4").Value = $dtZero
    $ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("G4").Value = $dtZero
    $ws.Range("H4").Value = $statusIgnored

    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBlob/$docCommit/e2e/$newDocMd", "", "", $newDocMd)
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBlob/$docCommit/$($locale.SheetName)/$xlfDisplay", "", "", $xlfDisplay)

    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBlob/$docCommit/e2e/$newCfgMd", "", "", $newCfgMd)
    $ws.Hyperlinks.Add($ws.Range("C3"), "$repoBlob/$docCommit/$($locale.SheetName)/$xlfDisplay", "", "", $xlfDisplay)

    $ws.Hyperlinks.Add($ws.Range("A4"), "$repoBlob/$cfgCommit/$localConfig", "", "", $localConfig)
}

Write-Host "Generated handoff report"
